$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.695.63'
$ws.Range('E2').Value = '  +3.39%  '
$ws.Range('D3').Value = '1.860.66'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '''230.70'
$ws.Range('E5').Value = '  +1.92%  '
$ws.Range('D6').Value = '''0.612'
$ws.Range('E6').Value = '  +2.89%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').Value = '''41.99'
$ws.Range('E8').Value = '  +9.10%  '
$ws.Range('E9').Value = '  +6.97%  '
$ws.Range('D10').Value = '''0.0696'
$ws.Range('E10').Value = '  +2.55%  '
$ws.Range('E11').Value = '  +2.84%  '
$ws.Range('D12').Value = '2.134.04'
$ws.Range('E12').Value = '  +2.74%  '
$ws.Range('D13').Value = '''11.51'
$ws.Range('E13').Value = '  +1.84%  '
$ws.Range('D14').Value = '1.870.18'
$ws.Range('E14').Value = '  +2.31%  '
$ws.Range('E15').Value = '  +6.85%  '
$ws.Range('D16').Value = '''4.77'
$ws.Range('E16').Value = '  +7.34%  '
$ws.Range('D17').Value = '35.644.33'
$ws.Range('E17').Value = '  +3.33%  '
$ws.Range('D18').Value = '''70.21'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').Value = '0.0₃0801'
$ws.Range('E19').Value = '  +3.09%  '
$ws.Range('D20').Value = '''246.44'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').Value = '''12.19'
$ws.Range('E21').Value = '  +8.16%  '
$ws.Range('D22').Value = '''4.77'
$ws.Range('E22').Value = '  +15.01%  '
$ws.Range('D24').Value = '''2.22'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = '''171.41'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('D26').Value = '''7.92'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('D27').Value = '''17.87'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').Value = '''0.123'
$ws.Range('E28').Value = '  +2.33%  '
$ws.Range('D29').Value = '''1.43'
$ws.Range('E29').Value = '  +15.86%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').Value = '3.367.07'
$ws.Range('E31').Value = '  +38.58%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '''0.0545'
$ws.Range('E32').Value = '  +4.69%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''3.93'
$ws.Range('E33').Value = '  +3.47%  '
$ws.Range('D34').Value = '''4.05'
$ws.Range('E34').Value = '  +3.91%  '
$ws.Range('D35').Value = '''1.90'
$ws.Range('E35').Value = '  +3.51%  '
$ws.Range('D36').Value = '''0.692'
$ws.Range('E36').Value = '  +6.16%  '
$ws.Range('E37').Value = '  +5.88%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '''1.09'
$ws.Range('E38').Value = '  +2.91%  '
$ws.Range('D39').Value = '''88.73'
$ws.Range('E39').Value = '  +8.55%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.344.48'
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.0196'
$ws.Range('E41').Value = '  +4.77%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '''1.03'
$ws.Range('E42').Value = '  +7.21%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '''1.29'
$ws.Range('E43').Value = '  +6.45%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '''15.12'
$ws.Range('E44').Value = '  +8.84%  '
$ws.Range('D45').Value = '''2.47'
$ws.Range('E45').Value = '  +1.47%  '
$ws.Range('D46').Value = '''2.83'
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('E47').Value = '  +3.07%  '
$ws.Range('D48').Value = '''6.10'
$ws.Range('E48').Value = '  +5.28%  '
$ws.Range('D49').Value = '2.031.03'
$ws.Range('E49').Value = '  +2.66%  '
$ws.Range('D50').Value = '''104.72'
$ws.Range('E50').Value = '  +2.21%  '
$ws.Range('E51').Value = '  +0.22%  '
